# Applies the odds updates described in the commit diff to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("O3").Value = 1.73
$ws.Range("P3").Value = 2

# Row 4
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3

# Row 5
$ws.Range("Q5").Value = 1.75
$ws.Range("R5").Value = 2.05

# Row 6
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 3.5
$ws.Range("K6").Value = 2.15
$ws.Range("L6").Value = 2.18
$ws.Range("N6").Value = 8.1
$ws.Range("T6").Value = 2.57
$ws.Range("U6").Value = 1.98
$ws.Range("V6").Value = 1.65
$ws.Range("X6").Value = 28
$ws.Range("Z6").Value = 90
$ws.Range("AD6").Value = 6.9
$ws.Range("AE6").Value = 18.5
$ws.Range("AH6").Value = 7
$ws.Range("AJ6").Value = 12
$ws.Range("AM6").Value = 1000
$ws.Range("AO6").Value = 28
$ws.Range("AP6").Value = 35
$ws.Range("AT6").Value = 2.55
$ws.Range("AU6").Value = 7.7
$ws.Range("AV6").Value = 75
$ws.Range("AX6").Value = 7.9
$ws.Range("AZ6").Value = 26

# Row 8
$ws.Range("Q8").Value = 2.05
$ws.Range("R8").Value = 1.75

# Row 11
$ws.Range("G11").Value = 1.9
$ws.Range("I11").Value = 3.7
$ws.Range("J11").Value = 2.6
$ws.Range("K11").Value = 2.2
$ws.Range("X11").Value = 9
$ws.Range("AA11").Value = 15
$ws.Range("AK11").Value = 34
$ws.Range("AL11").Value = 41
$ws.Range("AQ11").Value = 34
$ws.Range("BA11").Value = 101
